# Updated cryptos list on Mon Sep 25 23:27:52 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) text cells
# for the rows whose figures moved since the last scrape. All cells in
# these columns are stored as plain text (not numbers), so values that
# look numeric (e.g. "210.21", "19.40") must be forced back to text -
# otherwise Excel's COM layer will auto-coerce them to the Number type
# and silently drop significant trailing zeros / thousands-style dots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param([string]$Ref, [string]$Value)

    $cell = $ws.Range($Ref)

    # Only cells whose new value would be auto-parsed as a plain number
    # (single decimal point, no letters) need the text-format round trip;
    # values with two dots ("1.588.75"), a percent sign, spaces, etc. are
    # already unambiguous text and are left with their original (no-op)
    # cell style.
    $isPlainNumber = $Value -match '^[+-]?[0-9]+(\.[0-9]+)?$'

    if ($isPlainNumber) {
        $cell.NumberFormat = '@'
        $cell.Value = $Value
        # Drop the Text (@) number format we just applied so the cell's
        # style index matches the untouched cells around it.
        $cell.ClearFormats()
    } else {
        $cell.Value = $Value
    }
}

Set-TextCell 'D2' '26.318.11'
Set-TextCell 'E2' '  -0.93%  '
Set-TextCell 'D3' '1.588.75'
Set-TextCell 'E3' '  -0.27%  '
Set-TextCell 'E4' '  -0.51%  '
Set-TextCell 'D5' '210.21'
Set-TextCell 'E5' '  -0.03%  '
Set-TextCell 'E6' '  -0.36%  '
Set-TextCell 'E7' '  -0.46%  '
Set-TextCell 'E8' '  -0.57%  '
Set-TextCell 'D9' '0.246'
Set-TextCell 'E9' '  +0.04%  '
Set-TextCell 'D10' '19.40'
Set-TextCell 'E10' '  -0.55%  '
Set-TextCell 'E11' '  +0.14%  '
Set-TextCell 'E12' '  -0.29%  '
Set-TextCell 'E13' '  +0.74%  '
Set-TextCell 'D14' '1.584.28'
Set-TextCell 'E14' '  -1.37%  '
Set-TextCell 'D15' '0.519'
Set-TextCell 'E15' '  +0.15%  '
Set-TextCell 'D16' '64.36'
Set-TextCell 'E16' '  +0.02%  '
Set-TextCell 'E17' '  -0.91%  '
Set-TextCell 'E18' '  -1.05%  '
Set-TextCell 'E19' '  +5.87%  '
Set-TextCell 'D20' '210.97'
Set-TextCell 'E20' '  +1.63%  '
Set-TextCell 'E21' '  -0.42%  '
Set-TextCell 'E22' '  -0.10%  '
Set-TextCell 'D23' '8.95'
Set-TextCell 'E23' '  +0.52%  '
Set-TextCell 'E24' '  -3.00%  '
Set-TextCell 'D25' '144.54'
Set-TextCell 'E25' '  -0.32%  '
Set-TextCell 'E26' '  -0.49%  '
Set-TextCell 'E27' '  -0.54%  '
Set-TextCell 'E28' '  -0.45%  '
Set-TextCell 'D29' '15.25'
Set-TextCell 'E29' '  -0.19%  '
Set-TextCell 'E30' '  +0.41%  '
Set-TextCell 'D32' '3.22'
Set-TextCell 'E32' '  -0.69%  '
Set-TextCell 'E33' '  +1.67%  '
Set-TextCell 'D34' '1.317.12'
Set-TextCell 'E34' '  +2.92%  '
Set-TextCell 'E35' '  -1.49%  '
Set-TextCell 'D36' '0.609'
Set-TextCell 'E36' '  +1.76%  '
Set-TextCell 'E37' '  -0.46%  '
Set-TextCell 'E38' '  +0.37%  '
Set-TextCell 'E40' '  -1.45%  '
Set-TextCell 'E41' '  -0.46%  '
Set-TextCell 'E42' '  +3.92%  '
Set-TextCell 'E43' '  -0.31%  '
Set-TextCell 'E44' '  -1.00%  '
Set-TextCell 'D45' '62.30'
Set-TextCell 'E45' '  -0.22%  '
Set-TextCell 'D46' '1.725.09'
Set-TextCell 'E46' '  -0.27%  '
Set-TextCell 'D47' '87.53'
Set-TextCell 'E47' '  -1.95%  '
Set-TextCell 'E48' '  -5.16%  '
Set-TextCell 'E49' '  -1.19%  '
Set-TextCell 'D50' '0.0976'
Set-TextCell 'E50' '  -4.61%  '
Set-TextCell 'E51' '  -0.42%  '
